$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. PersonalAutoData: selection moves from I1 to C2
# ---------------------------------------------------------------------------
$wsPersonal = $wb.Worksheets.Item("PersonalAutoData")
[void]$wsPersonal.Range("C2").Select()

# ---------------------------------------------------------------------------
# 2. VehicleData: selection moves from L11 to the A1:C3 block
# ---------------------------------------------------------------------------
$wsVehicle = $wb.Worksheets.Item("VehicleData")
[void]$wsVehicle.Range("A1:C3").Select()

# ---------------------------------------------------------------------------
# 3. Add the new EndorsementData sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEndorsement = $wb.Worksheets.Add($null, $lastSheet)
$wsEndorsement.Name = "EndorsementData"

# Write the cell values first (in the same order the original strings were
# first introduced, so new shared-string entries line up: Description,
# Reason, Sub_Type, Type, Effective_Date, Endorsement, Policy Correction,
# Add Vehicle).
$wsEndorsement.Range("E1").Value = "Description"
$wsEndorsement.Range("C1").Value = "Reason"
$wsEndorsement.Range("B1").Value = "Sub_Type"
$wsEndorsement.Range("A1").Value = "Type"
$wsEndorsement.Range("D1").Value = "Effective_Date"

$wsEndorsement.Range("A2").Value = "Endorsement"
$wsEndorsement.Range("B2").Value = "Endorsement"
$wsEndorsement.Range("C2").Value = "Policy Correction"
$wsEndorsement.Range("E2").Value = "Add Vehicle"

# Apply the same formatting used on VehicleData: header row style, the
# quote-prefixed style used in VehicleData's first column (applied here to
# A2:A3 and B2), and the plain bordered style everywhere else.
$wsVehicle.Range("A1:E1").Copy()
[void]$wsEndorsement.Range("A1:E1").PasteSpecial(-4122)

$wsVehicle.Range("A2").Copy()
[void]$wsEndorsement.Range("A2").PasteSpecial(-4122)
[void]$wsEndorsement.Range("A3").PasteSpecial(-4122)
[void]$wsEndorsement.Range("B2").PasteSpecial(-4122)

$wsVehicle.Range("B2").Copy()
[void]$wsEndorsement.Range("C2").PasteSpecial(-4122)
[void]$wsEndorsement.Range("D2").PasteSpecial(-4122)
[void]$wsEndorsement.Range("E2").PasteSpecial(-4122)
[void]$wsEndorsement.Range("B3").PasteSpecial(-4122)
[void]$wsEndorsement.Range("C3").PasteSpecial(-4122)
[void]$wsEndorsement.Range("D3").PasteSpecial(-4122)
[void]$wsEndorsement.Range("E3").PasteSpecial(-4122)

# Column widths (character-width units; the host stores width = ColumnWidth
# + 5/6, quantized to 1/6 steps, so these values land as close as possible
# to the authored 15.6640625 / 20 / 27.5546875 / 18.109375 / 18.88671875).
$wsEndorsement.Columns.Item(1).ColumnWidth = 14.833333333333334
$wsEndorsement.Columns.Item(2).ColumnWidth = 19.166666666666668
$wsEndorsement.Columns.Item(3).ColumnWidth = 26.666666666666668
$wsEndorsement.Columns.Item(4).ColumnWidth = 17.333333333333332
$wsEndorsement.Columns.Item(5).ColumnWidth = 18.0

# Final selection on the new sheet (it also becomes the active sheet/tab,
# matching the authored file).
[void]$wsEndorsement.Range("E3").Select()
